# Apply the edits described by the diff:
#  - Insert a new "2018" data row above the existing 2019 row, shifting the
#    rest of the table (2019-2024) down by one row.
#  - Move the existing cell comment from C7 to C8 so it still annotates the
#    2024 Qtr2 value that shifted down.
#  - Update the active cell selection.
#  - Update the "Generated on" timestamp in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert the new row for 2018 data, pushing 2019-2024 down one row ---
$ws.Rows.Item(2).Insert()

# Copy the number/cell formatting from the row below (now row 3, the old
# 2019 row) onto the newly inserted row 2 so it matches the rest of the
# table instead of inheriting the header's style.
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = 2018
$ws.Range("B2").Value = -0.8
$ws.Range("C2").Value = -0.1
$ws.Range("D2").Value = 2.2
$ws.Range("E2").Value = -0.2

# --- Move the comment that was anchored on C7 down to C8 ---
$comment = $ws.Range("C7").Comment
$commentText = $comment.Text()
$comment.Delete()
$ws.Range("C8").AddComment($commentText)

# --- Update the selected cell shown in the saved view ---
$ws.Range("F18").Select()

# --- Update the footer generation timestamp ---
$ws.PageSetup.RightFooter = "Generated on: November 5, 2024 (12:53:05 AM)"
